# "added new window for the items"
# Append a new reservation record (row 7) to the "Order Data" sheet:
#   Reservation ID = 6, Number of Guests = 3,
#   Date Of Reservation = 2024-04-24 (serial 45406), Customer ID = 4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order Data")

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 45406
$ws.Range("C7").NumberFormat = "yyyy-MM-dd"
$ws.Range("D7").Value = 4
